$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-28 Sunday" "2025-12-29 Monday"

Replace-Text "36÷6=" "74÷9="
Replace-Text "40÷3=" "40÷9="
Replace-Text "15÷9=" "93÷9="
Replace-Text "68÷2=" "43÷4="
Replace-Text "26÷3=" "95÷9="
Replace-Text "46÷4=" "33÷6="
Replace-Text "94÷9=" "33÷4="
Replace-Text "47÷4=" "71÷2="
Replace-Text "37÷2=" "27÷9="
Replace-Text "39÷9=" "59÷6="
Replace-Text "56÷7=" "44÷2="
Replace-Text "28÷4=" "78÷3="
Replace-Text "72÷9=" "97÷5="
Replace-Text "50÷3=" "96÷2="
Replace-Text "79÷5=" "27÷5="
Replace-Text "19÷9=" "99÷4="
Replace-Text "29÷6=" "86÷9="
Replace-Text "83÷4=" "84÷4="
Replace-Text "72÷6=" "54÷7="
Replace-Text "75÷6=" "37÷5="
Replace-Text "68÷5=" "79÷8="
Replace-Text "70÷4=" "37÷3="
Replace-Text "99÷2=" "39÷7="
Replace-Text "16÷3=" "13÷2="
Replace-Text "35÷9=" "11÷3="
